$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.531.76'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +4.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.844.22'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +3.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.031'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +3.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.55'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +4.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.029'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +2.85%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4377'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.29%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3738'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +3.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07397'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8770'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.77%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.42'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.866.25'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +6.24%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.488'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +4.71%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.694'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07164'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.01%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '82.73'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.62%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.033'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009036'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.68%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.028'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.67%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.36'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.547.96'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.232'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +2.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.22'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +3.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.071.21'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +4.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '157.09'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +3.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.930'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.66%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.71'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.255'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.939'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.60%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '116.34'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.71%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09075'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.208'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +7.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7631'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +5.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.485'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.93%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.882'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +5.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.031'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.149'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +6.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01970'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.58%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05256'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.88%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5180'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.30%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.771'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +6.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1662'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '6.615'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +4.59%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.503'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +5.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.15'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.33%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.61'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +3.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.032'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +3.20%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +4.41%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4640'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +3.94%  '
$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06334'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.53%  '
$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.886'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +9.58%  '
